$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Adatok")

# Add a new row (row 7) with data about a captured image upload option ("Probakép")
$ws.Range("A7").Value = "Ágynemű"
$ws.Range("B7").Value = "Probakép"
$ws.Range("C7").Value = "64-110"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "igen"
$ws.Range("F7").Value = "Saját"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "2023-05-13"
$ws.Range("G7").Style = "Normal"
